$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '67.430.63'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.97%  '
$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '3.594.68'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.90%  '
$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '578.86'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.80%  '
$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '193.05'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '3.591.61'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.88%  '
$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.618'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.678'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.15%  '
$ws.Range("E11").Value = '  -1.28%  '
$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '54.55'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.86%  '
$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '0.0000275'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '9.95'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.48%  '
$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '4.159.40'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.96%  '
$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '3.587.72'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.00%  '
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '12.34'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '18.43'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.41%  '
$ws.Range("B20").Value = 'WrappedBTC'
$ws.Range("C20").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '67.359.22'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.75%  '
$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '1.08'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -3.25%  '
$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '402.66'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '13.34'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +20.49%  '
$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '4.21'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.81%  '
$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '85.80'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.26%  '
$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '2.92'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.01%  '
$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '12.60'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '6.10'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.60%  '
$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '3.80'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.61%  '
$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '8.09'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +11.11%  '
$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '9.18'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.07%  '
$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '31.36'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.59%  '
$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '668.83'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +10.44%  '
$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '12.24'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '64.15'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.87%  '
$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '42.87'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.49%  '
$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.427'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +8.60%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.0₃0788'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +2.62%  '
$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '2.95'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +16.61%  '
$ws.Range("E42").Value = '  +8.37%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '0.134'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.02%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '3.155.11'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +13.70%  '
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.0419'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.80%  '
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '3.14'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '8.80'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '142.59'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("E51").Value = '  -3.49%  '
